# Update "想去人数" (interested-count) values in column F across the
# "展览", "演出" and "全部类型" worksheets, reflecting refreshed scrape data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1375
$ws.Range("F3").Value = 2013
$ws.Range("F4").Value = 910
$ws.Range("F5").Value = 117
$ws.Range("F7").Value = 693
$ws.Range("F8").Value = 121
$ws.Range("F10").Value = 40
$ws.Range("F11").Value = 2531
$ws.Range("F12").Value = 1612
$ws.Range("F13").Value = 1578
$ws.Range("F15").Value = 255
$ws.Range("F16").Value = 636
$ws.Range("F17").Value = 810
$ws.Range("F18").Value = 97
$ws.Range("F19").Value = 322
$ws.Range("F20").Value = 1095
$ws.Range("F22").Value = 36
$ws.Range("F23").Value = 531
$ws.Range("F24").Value = 5407
$ws.Range("F25").Value = 227
$ws.Range("F26").Value = 762
$ws.Range("F27").Value = 97
$ws.Range("F29").Value = 145
$ws.Range("F30").Value = 238
$ws.Range("F31").Value = 227
$ws.Range("F32").Value = 42
$ws.Range("F33").Value = 1054
$ws.Range("F34").Value = 771
$ws.Range("F36").Value = 61
$ws.Range("F38").Value = 410
$ws.Range("F39").Value = 1137
$ws.Range("F40").Value = 141
$ws.Range("F42").Value = 182
$ws.Range("F43").Value = 133
$ws.Range("F44").Value = 81

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 789
$ws.Range("F5").Value = 432
$ws.Range("F6").Value = 12
$ws.Range("F7").Value = 10

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1375
$ws.Range("F4").Value = 2013
$ws.Range("F5").Value = 910
$ws.Range("F7").Value = 117
$ws.Range("F9").Value = 693
$ws.Range("F10").Value = 121
$ws.Range("F11").Value = 12
$ws.Range("F12").Value = 10
$ws.Range("F14").Value = 40
$ws.Range("F15").Value = 2531
$ws.Range("F16").Value = 1612
$ws.Range("F17").Value = 1578
$ws.Range("F19").Value = 255
$ws.Range("F20").Value = 636
$ws.Range("F22").Value = 810
$ws.Range("F23").Value = 97
$ws.Range("F24").Value = 322
$ws.Range("F25").Value = 1095
$ws.Range("F26").Value = 36
$ws.Range("F27").Value = 531
$ws.Range("F28").Value = 5407
$ws.Range("F29").Value = 227
$ws.Range("F30").Value = 762
$ws.Range("F31").Value = 97
$ws.Range("F33").Value = 145
$ws.Range("F34").Value = 238
$ws.Range("F35").Value = 227
$ws.Range("F36").Value = 42
$ws.Range("F37").Value = 1054
$ws.Range("F38").Value = 771
$ws.Range("F39").Value = 61
$ws.Range("F40").Value = 410
$ws.Range("F41").Value = 1137
$ws.Range("F42").Value = 141
$ws.Range("F44").Value = 182
$ws.Range("F45").Value = 133
$ws.Range("F46").Value = 81
